# Add a "booth number" (攤位編號) column into the club import sample sheet.
# This inserts a new column D (shifting the existing 攤位負責人1-5 columns
# from D:H to E:I), fills in the header/value for the new column, restores
# the narrower "id-style" column width used by the other id columns, and
# re-points the conditional-formatting range that covered the shifted
# columns so it keeps covering the same (now shifted) header columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column D; this shifts D:H -> E:I,
# along with their data, styles, and column widths.
$ws.Columns("D:D").Insert()

# Populate the new column's header and sample value.
$ws.Range("D1").Value = "攤位編號"
$ws.Range("D2").Value = "D15"

# The newly inserted column should look like the other "id" columns
# (A-C: width 10.25) rather than the wider bestFit columns it pushed
# to the right (E-I: width 13.875).
$ws.Columns("D:D").ColumnWidth = 9.5

# The conditional formatting rule that used to cover D1:H1048576 now
# needs to cover the shifted range E1:I1048576 (the new column D is
# intentionally left out of this duplicate-value highlighting rule).
$oldRuleRange = $ws.Range("D1:H1048576")
$fcs = $oldRuleRange.FormatConditions
$fc = $fcs.Item(1)
$fc.ModifyAppliesToRange($ws.Range("E1:I1048576"))

# Match the saved selection/active cell.
$ws.Range("G7").Select() | Out-Null
